$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.996.89"
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").Value = "2.654.19"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'523.93"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").Value = "'144.34"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("D9").Value = "'7.02"
$ws.Range("E9").Value = "  +9.00%  "
$ws.Range("E10").Value = "  -2.64%  "
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "'0.130"
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("D13").Value = "3.122.99"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "59.039.93"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "'21.04"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.674.44"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000136"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "'339.42"
$ws.Range("E18").Value = "  -3.24%  "
$ws.Range("D19").Value = "'4.36"
$ws.Range("E19").Value = "  -4.22%  "
$ws.Range("D20").Value = "'10.35"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "'0.996"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'63.98"
$ws.Range("E23").Value = "  +2.07%  "
$ws.Range("D24").Value = "'0.418"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "0.0₃0801"
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("D28").Value = "'7.08"
$ws.Range("E28").Value = "  -2.48%  "
$ws.Range("D29").Value = "'6.67"
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "'18.86"
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("D33").Value = "'149.75"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  -3.68%  "
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("D36").Value = "'0.890"
$ws.Range("E36").Value = "  -6.29%  "
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").Value = "'36.80"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("E39").Value = "  -5.99%  "
$ws.Range("E40").Value = "  -3.09%  "
$ws.Range("D41").Value = "'0.618"
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").Value = "'19.96"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "'276.17"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("D48").Value = "2.034.51"
$ws.Range("E48").Value = "  -4.76%  "
$ws.Range("D49").Value = "'4.75"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").Value = "'0.0228"
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("D51").Value = "'18.80"
$ws.Range("E51").Value = "  -1.61%  "

Write-Output "Applied cryptos update."
